$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from SCD0287 to SCD0018
$ws.Name = "SCD0018"

# Update B2 TC_ID value, copying style (font/alignment) from C2 (TEST_SCENARIO_DESC cell)
$ws.Range("B2").Value = "SCD0018-010"
$ws.Range("C2").Copy()
$ws.Range("B2").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B2").Value = "SCD0018-010"

# Adjust column B width
$ws.Columns.Item(2).ColumnWidth = 12.28515625

# Update the active selection to C2
$ws.Range("C2").Select()
